# floridaCityCoordinates.xlsx edit
#
# Source diff summary:
#  - Row 891 ("Whitfield", lat 30.88 / long -87.06 -- a duplicate Whitfield
#    entry; shared-string index 887 already appears on row 890 with
#    different coordinates) is removed from Sheet1. Every row below it
#    shifts up by one, so the sheet goes from 920 data rows (+1 header) to
#    919 data rows (+1 header), and the used range shrinks from A1:C921 to
#    A1:C920.
#  - The selected cell moves from H7 to H16.
#  - An AutoFilter is turned on over A1:C920, which Excel always backs with
#    a hidden, sheet-scoped "_FilterDatabase" defined name.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the duplicated "Whitfield" row -------------------------------
[void]$ws.Rows("891:891").Delete()

# --- Turn on AutoFilter across the (now smaller) data range --------------
[void]$ws.Range("A1:C920").AutoFilter()

# Excel always persists the AutoFilter's range as a hidden, sheet-scoped
# _FilterDatabase defined name; add it explicitly so it round-trips on save.
$filterDatabaseName = $ws.Names.Add("_xlnm._FilterDatabase", "=Sheet1!`$A`$1:`$C`$920")
$filterDatabaseName.Visible = $false

# --- Move the active selection to H16 -------------------------------------
[void]$ws.Range("H16").Select()
